$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.765.79"
$ws.Range("E2").Value = "  +2.24%  "
$ws.Range("D3").Value = "3.735.38"
$ws.Range("E3").Value = "  +19.25%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'616.11"
$ws.Range("E5").Value = "  +6.67%  "
$ws.Range("D6").Value = "'177.34"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("D7").Value = "3.731.08"
$ws.Range("E7").Value = "  +19.14%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  +4.69%  "
$ws.Range("D10").Value = "'0.168"
$ws.Range("E10").Value = "  +9.85%  "
$ws.Range("D11").Value = "'6.40"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").Value = "'0.503"
$ws.Range("E12").Value = "  +7.41%  "
$ws.Range("D13").Value = "'40.90"
$ws.Range("E13").Value = "  +11.12%  "
$ws.Range("D14").Value = "'0.0000256"
$ws.Range("E14").Value = "  +5.57%  "
$ws.Range("D15").Value = "4.362.01"
$ws.Range("E15").Value = "  +19.41%  "
$ws.Range("D16").Value = "3.739.15"
$ws.Range("E16").Value = "  +19.43%  "
$ws.Range("D17").Value = "69.870.18"
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").Value = "'7.60"
$ws.Range("E19").Value = "  +6.80%  "
$ws.Range("D20").Value = "'517.13"
$ws.Range("E20").Value = "  +6.40%  "
$ws.Range("D21").Value = "'16.72"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").Value = "'9.40"
$ws.Range("E22").Value = "  +20.41%  "
$ws.Range("D23").Value = "'0.737"
$ws.Range("E23").Value = "  +5.85%  "
$ws.Range("D24").Value = "'88.78"
$ws.Range("E24").Value = "  +5.87%  "
$ws.Range("D25").Value = "'2.49"
$ws.Range("E25").Value = "  +6.67%  "
$ws.Range("D26").Value = "'13.53"
$ws.Range("E26").Value = "  +4.22%  "
$ws.Range("D27").Value = "'10.91"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "'0.0000125"
$ws.Range("E29").Value = "  +31.99%  "
$ws.Range("D30").Value = "'2.50"
$ws.Range("E30").Value = "  +5.90%  "
$ws.Range("D31").Value = "'2.86"
$ws.Range("E31").Value = "  +8.39%  "
$ws.Range("D32").Value = "'7.85"
$ws.Range("E32").Value = "  -3.59%  "
$ws.Range("D33").Value = "'31.55"
$ws.Range("E33").Value = "  +11.93%  "
$ws.Range("D34").Value = "'0.115"
$ws.Range("E34").Value = "  +2.06%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "'6.21"
$ws.Range("E36").Value = "  +9.82%  "
$ws.Range("E37").Value = "  +8.66%  "
$ws.Range("E38").Value = "  +5.25%  "
$ws.Range("D39").Value = "'2.19"
$ws.Range("E39").Value = "  +6.51%  "
$ws.Range("E40").Value = "  +6.27%  "
$ws.Range("D41").Value = "'51.47"
$ws.Range("E41").Value = "  +4.73%  "
$ws.Range("E42").Value = "  +5.58%  "
$ws.Range("D43").Value = "'44.57"
$ws.Range("E43").Value = "  -8.64%  "
$ws.Range("D44").Value = "'423.15"
$ws.Range("E44").Value = "  +6.15%  "
$ws.Range("D45").Value = "3.075.51"
$ws.Range("E45").Value = "  +9.97%  "
$ws.Range("D46").Value = "'2.73"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").Value = "'0.0366"
$ws.Range("E47").Value = "  +4.93%  "
$ws.Range("E48").Value = "  +2.07%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'136.81"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.52"
$ws.Range("E50").Value = "  +7.02%  "
$ws.Range("E51").Value = "  -0.01%  "
